$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Move the section-header labels from column C to column B (rows 16/25)
# ---------------------------------------------------------------------
$v16 = $ws.Range("C16").Value()
$ws.Range("C16").ClearContents()
$ws.Range("B16").Value = $v16

$v25 = $ws.Range("C25").Value()
$ws.Range("C25").ClearContents()
$ws.Range("B25").Value = $v25

# ---------------------------------------------------------------------
# Helper color value for the "FF9900" highlight slice (BGR OLE color)
# 0xFF9900 -> R=255 G=153 B=0 -> R + G*256 + B*65536
# ---------------------------------------------------------------------
$highlightColor = 255 + (153 * 256) + (0 * 65536)

# ---------------------------------------------------------------------
# Chart 1 : "Lieu du stage"  (Worksheet!$B$10, $D$10:$D$14, $E$10:$E$14)
# ---------------------------------------------------------------------
$chartObj1 = $ws.ChartObjects().Add(450, 20, 430, 270)
$chart1 = $chartObj1.Chart
$chart1.ChartType = 5
$ser1 = $chart1.SeriesCollection().NewSeries()
$ser1.Formula = "=SERIES(Worksheet!`$B`$10,Worksheet!`$D`$10:`$D`$14,Worksheet!`$E`$10:`$E`$14,1)"
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Lieu du stage"
$chart1.HasLegend = $true
$chart1.Legend.Position = -4152
$ser1.VaryByCategories = $true
$ser1.HasDataLabels = $true
$dl1 = $ser1.DataLabels()
$dl1.ShowValue = $true
$dl1.ShowPercentage = $true
$dl1.ShowLegendKey = $false
$dl1.ShowCategoryName = $false
$dl1.ShowSeriesName = $false
$dl1.ShowBubbleSize = $false
$pt1 = $ser1.Points(4)
$pt1.Format.Fill.ForeColor.RGB = $highlightColor
$chartObj1.Name = "Chart 1"

# ---------------------------------------------------------------------
# Chart 2 : "Contenu du stage" (Worksheet!$B$16, $D$16:$D$23, $E$16:$E$23)
# ---------------------------------------------------------------------
$chartObj2 = $ws.ChartObjects().Add(450, 320, 430, 270)
$chart2 = $chartObj2.Chart
$chart2.ChartType = 5
$ser2 = $chart2.SeriesCollection().NewSeries()
$ser2.Formula = "=SERIES(Worksheet!`$B`$16,Worksheet!`$D`$16:`$D`$23,Worksheet!`$E`$16:`$E`$23,1)"
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Contenu du stage"
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152
$ser2.VaryByCategories = $true
$ser2.HasDataLabels = $true
$dl2 = $ser2.DataLabels()
$dl2.ShowValue = $true
$dl2.ShowPercentage = $true
$dl2.ShowLegendKey = $false
$dl2.ShowCategoryName = $false
$dl2.ShowSeriesName = $false
$dl2.ShowBubbleSize = $false
$pt2 = $ser2.Points(4)
$pt2.Format.Fill.ForeColor.RGB = $highlightColor
$chartObj2.Name = "Chart 2"

# ---------------------------------------------------------------------
# Chart 3 : "Type du stage" (Worksheet!$B$25, $D$25:$D$28, $E$25:$E$28)
# ---------------------------------------------------------------------
$chartObj3 = $ws.ChartObjects().Add(450, 620, 430, 270)
$chart3 = $chartObj3.Chart
$chart3.ChartType = 5
$ser3 = $chart3.SeriesCollection().NewSeries()
$ser3.Formula = "=SERIES(Worksheet!`$B`$25,Worksheet!`$D`$25:`$D`$28,Worksheet!`$E`$25:`$E`$28,1)"
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Type du stage"
$chart3.HasLegend = $true
$chart3.Legend.Position = -4152
$ser3.VaryByCategories = $true
$ser3.HasDataLabels = $true
$dl3 = $ser3.DataLabels()
$dl3.ShowValue = $true
$dl3.ShowPercentage = $true
$dl3.ShowLegendKey = $false
$dl3.ShowCategoryName = $false
$dl3.ShowSeriesName = $false
$dl3.ShowBubbleSize = $false
$pt3 = $ser3.Points(4)
$pt3.Format.Fill.ForeColor.RGB = $highlightColor
$chartObj3.Name = "Chart 3"

Write-Host "Charts added"
